$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questions")
$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$text = @'
questions = [
    {
        "title": "You are working on a project and have made some changes to the code. You plan to make changes to some other areas of the code as well. When should you commit changes to the code?",
        "ques_type": 2,
        "options": [
            "After you are finished with each feature/task.",
            "After changing each file.",
            "At the end, after you have made all the changes.",
            "After the project is complete."
        ],
        "score": "After you are finished with each feature/task."
    },
    {
        "title": "You are planning to use an existing open-source project as the basis for your own idea. You intend to improve the existing project and, in the future, propose the changes to the original project owner. What should be your first steps in accomplishing this?",
        "ques_type": 2,
        "options": [
            "Fork the existing open-source project and make changes to it.",
            "Clone the existing open-source project and push changes to it.",
            "Download the existing open-source project and make changes to it.",
            "Create a new branch on the existing open-source project, make changes to it, and merge the changes later."
        ],
        "score": "Fork the existing open-source project and make changes to it."
    },
    {
        "title": "In order to host a repository, GitHub Pages requires a specific file to be present in the repository. What is the name of this file?",
        "ques_type": 2,
        "options": [
            "main.js",
            "index.js",
            "start.html",
            "index.html"
        ],
        "score": "index.html"
    },
    {
        "title": "True or false: The GitHub Issues feature is used to run tests on code in order to identify any issues in the code.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "False"
    }
]
'@
$ws.Range("A1").Value = $text
$ws.Rows(1).AutoFit()
